# Weekly refresh: prepend two new price observations (rows 409-410) to the
# Kiwi / Macroferia Regional de Talca log, pushing the existing rows 409-426
# down to 411-428.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 409; this shifts rows
# 409:426 down to 411:428, preserving their formatting/styles.
$ws.Rows("409:410").Insert()

# --- New row 409 -----------------------------------------------------
$ws.Range("A409").Value2 = 5
$ws.Range("B409").Value2 = "Macroferia Regional de Talca"
$ws.Range("C409").Value2 = "Maule"
$ws.Range("D409").Value2 = "2023-05-29"
$ws.Range("E409").Value2 = 7
$ws.Range("F409").Value2 = "Fruta"
$ws.Range("G409").Value2 = 100101
$ws.Range("H409").Value2 = "Berries"
$ws.Range("I409").Value2 = 100101007
$ws.Range("J409").Value2 = "Kiwi"
$ws.Range("K409").Value2 = "Hayward"
$ws.Range("L409").Value2 = "Especial"
$ws.Range("M409").Value2 = 220
$ws.Range("N409").Value2 = 10000
$ws.Range("O409").Value2 = 10000
$ws.Range("P409").Value2 = 10000
$ws.Range("Q409").Value2 = "$/bandeja 18 kilos"
$ws.Range("R409").Value2 = "Provincia de Curicó"
$ws.Range("S409").Value2 = 556
$ws.Range("T409").Value2 = 18

# --- New row 410 -----------------------------------------------------
$ws.Range("A410").Value2 = 5
$ws.Range("B410").Value2 = "Macroferia Regional de Talca"
$ws.Range("C410").Value2 = "Maule"
$ws.Range("D410").Value2 = "2023-05-29"
$ws.Range("E410").Value2 = 7
$ws.Range("F410").Value2 = "Fruta"
$ws.Range("G410").Value2 = 100101
$ws.Range("H410").Value2 = "Berries"
$ws.Range("I410").Value2 = 100101007
$ws.Range("J410").Value2 = "Kiwi"
$ws.Range("K410").Value2 = "Hayward"
$ws.Range("L410").Value2 = "Primera"
$ws.Range("M410").Value2 = 640
$ws.Range("N410").Value2 = 8000
$ws.Range("O410").Value2 = 9000
$ws.Range("P410").Value2 = 8719
$ws.Range("Q410").Value2 = "$/bandeja 18 kilos"
$ws.Range("R410").Value2 = "Provincia de Curicó"
$ws.Range("S410").Value2 = 484
$ws.Range("T410").Value2 = 18
